$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text formatting (many values look
# numeric, e.g. '83.04' or '0.9989', and would otherwise be auto-converted to
# numbers by Excel's value parser, dropping formatting like trailing zeros).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.409.59'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.846.70'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '240.67'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6315'
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07557'
$ws.Range('E8').Value = '  -0.65%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2954'
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '24.50'
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.852.88'
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6845'
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.00001005'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '83.04'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.123'
$ws.Range('E17').Value = '  -2.56%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '29.437.19'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '228.33'
$ws.Range('E19').Value = '  -2.80%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.47'
$ws.Range('E20').Value = '  -0.99%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9997'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.537'
$ws.Range('E22').Value = '  -1.84%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '157.02'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1392'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.374'
$ws.Range('E26').Value = '  -1.42%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.67'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.469'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.05730'
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.264'
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.120'
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.024'
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.846'
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.155'
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7113'
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.588'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.250.40'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01808'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.774'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9102'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.176'
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '101.69'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '66.08'
$ws.Range('E44').Value = '  -3.14%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00000000119'
$ws.Range('E45').Value = '  +2.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.090'
$ws.Range('E46').Value = '  -3.97%  '
$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4019'
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.090'
$ws.Range('E48').Value = '  -1.13%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.681'
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.1120'
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05727'
$ws.Range('E51').Value = '  -0.61%  '
